$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PMTestData")

# Column A labels first (so shared-string indices are allocated in the
# same order Excel would use: A44, A45, then B44, B45).
$ws.Cells.Item(44, 1).Value = "test_clearSCA_funcKey"
$ws.Cells.Item(45, 1).Value = "test_clearSCABR_funcKey"

$ws.Cells.Item(44, 2).Value = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function SCA --display-text "SCA" --line 1 --line-dir 10000,extension_key -i --dir 10001 --function SCA --display-text "SCA" --key 1 --line-dir 10000,10000,10001,extension_key -e -d 10000 --key 1,extension_key -e -d 10001 --line 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'
$ws.Cells.Item(45, 2).Value = 'number_initiate -number 10000..10001 -numbertype ex,extension -i -d 10000..10001 -l 1 --csp 0,ip_extension -i -d 10000..10001,extension_key -i --dir 10000 --function SCABR --display-text "SCABR" --line 1 --line-dir 10000,extension_key -i --dir 10001 --function SCABR --display-text "SCABR" --key 1 --line-dir 10000,10000,10001,extension_key -e -d 10000 --key 1,extension_key -e -d 10001 --line 1,ip_extension -e -d 10000..10001,extension -e -d 10000..10001,number_end -number 10000..10001 -numbertype ex'

$ws.Cells.Item(44, 3).Value = "Y"
$ws.Cells.Item(45, 3).Value = "Y"

# Match the wrap-text formatting used by the rest of column B.
$ws.Range("B44:B45").WrapText = $true

# These two rows hold long wrapped strings, so they get the same taller
# row height as the other multi-line rows above them.
$ws.Rows.Item(44).RowHeight = 116
$ws.Rows.Item(45).RowHeight = 116

# Update the view: scrolled down to show the new rows, with B45 selected.
$ws.Application.ActiveWindow.ScrollRow = 44
$ws.Range("B45").Select()
